# Convert the field-code runs ( fldChar begin/end + instrText runs ) that
# build the "m: ('<html...>').fromHTMLBodyString()" field into plain literal
# text runs using "{" and "}" instead of real field delimiters, exactly as
# TokenIteratorFieldRewriterSplit expects its tokens to look.

$d = $word.ActiveDocument

# Locate the field and the paragraph that contains it (defensive: don't
# hard-code a paragraph index, find it from the field's own position).
$f = $d.Fields.Item(1)
$codeStart = $f.Code.Start

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($codeStart -ge $p.Range.Start -and $codeStart -lt $p.Range.End) {
        $target = $p
        break
    }
}

# Rebuild the paragraph as literal-text runs (".fromHTMLBodyString()" field
# becomes "{ m: (' ... ' ).fromHTMLBodyString() }" in plain <w:t> runs),
# keeping the bookmark that sat in the middle of the field code untouched.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> ('</w:t></w:r><w:r><w:t>&lt;p&gt;First paragraph.&lt;br /&gt;\nSecond paragraph.&lt;/p&gt;\n\n&lt;p&gt;Third paragraph.&lt;/p&gt;\n\n&lt;p&gt;Last paragraph.&lt;/p&gt;</w:t></w:r><w:r><w:t>'</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>).from</w:t></w:r><w:r><w:t>HTML</w:t></w:r><w:r><w:t>Body</w:t></w:r><w:r><w:t>String()</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>
'@

$null = $target.Range.InsertXML($xml)
